$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "1;nitroglicerin,aspirin,atenolol,propranolol,rosuvastatin"
$ws.Range("A3").Value = "2;aspirin, nitroglicerin, propranolol"

$ws.Range("A3").Select() | Out-Null
